{"js": "const replacements = [\n  { oldText: \"2024-09-09 Monday\", newText: \"2024-09-10 Tuesday\" },\n  { oldText: \"375\u00f78=46, 7\", newText: \"165\u00f75=33, 0\" },\n  { oldText: \"621\u00f72=310, 1\", newText: \"111\u00f73=37, 0\" },\n  { oldText: \"973\u00f73=324, 1\", newText: \"155\u00f79=17, 2\" },\n  { oldText: \"311\u00f79=34, 5\", newText: \"848\u00f79=94, 2\" },\n  { oldText: \"939\u00f75=187, 4\", newText: \"514\u00f72=257, 0\" },\n  { oldText: \"573\u00f75=114, 3\", newText: \"180\u00f75=36, 0\" },\n  { oldText: \"564\u00f78=70, 4\", newText: \"123\u00f77=17, 4\" },\n  { oldText: \"797\u00f75=159, 2\", newText: \"672\u00f78=84, 0\" },\n  { oldText: \"954\u00f73=318, 0\", newText: \"760\u00f78=95, 0\" },\n  { oldText: \"230\u00f76=38, 2\", newText: \"467\u00f78=58, 3\" },\n  { oldText: \"148\u00f76=24, 4\", newText: \"411\u00f77=58, 5\" },\n  { oldText: \"290\u00f75=58, 0\", newText: \"979\u00f78=122, 3\" },\n  { oldText: \"340\u00f74=85, 0\", newText: \"760\u00f76=126, 4\" },\n  { oldText: \"913\u00f76=152, 1\", newText: \"816\u00f74=204, 0\" },\n  { oldText: \"658\u00f78=82, 2\", newText: \"507\u00f76=84, 3\" },\n  { oldText: \"925\u00f75=185, 0\", newText: \"816\u00f77=116, 4\" },\n  { oldText: \"286\u00f76=47, 4\", newText: \"958\u00f79=106, 4\" },\n  { oldText: \"267\u00f72=133, 1\", newText: \"335\u00f76=55, 5\" },\n  { oldText: \"902\u00f79=100, 2\", newText: \"341\u00f75=68, 1\" },\n  { oldText: \"380\u00f78=47, 4\", newText: \"144\u00f75=28, 4\" },\n  { oldText: \"807\u00f79=89, 6\", newText: \"994\u00f79=110, 4\" },\n  { oldText: \"294\u00f73=98, 0\", newText: \"188\u00f76=31, 2\" },\n  { oldText: \"706\u00f72=353, 0\", newText: \"647\u00f72=323, 1\" },\n  { oldText: \"980\u00f73=326, 2\", newText: \"134\u00f72=67, 0\" },\n  { oldText: \"732\u00f78=91, 4\", newText: \"372\u00f77=53, 1\" },\n];\n\nfor (const { oldText, newText } of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-09-09 Monday\"; New = \"2024-09-10 Tuesday\" }\n    @{ Old = \"375\u00f78=46, 7\"; New = \"165\u00f75=33, 0\" }\n    @{ Old = \"621\u00f72=310, 1\"; New = \"111\u00f73=37, 0\" }\n    @{ Old = \"973\u00f73=324, 1\"; New = \"155\u00f79=17, 2\" }\n    @{ Old = \"311\u00f79=34, 5\"; New = \"848\u00f79=94, 2\" }\n    @{ Old = \"939\u00f75=187, 4\"; New = \"514\u00f72=257, 0\" }\n    @{ Old = \"573\u00f75=114, 3\"; New = \"180\u00f75=36, 0\" }\n    @{ Old = \"564\u00f78=70, 4\"; New = \"123\u00f77=17, 4\" }\n    @{ Old = \"797\u00f75=159, 2\"; New = \"672\u00f78=84, 0\" }\n    @{ Old = \"954\u00f73=318, 0\"; New = \"760\u00f78=95, 0\" }\n    @{ Old = \"230\u00f76=38, 2\"; New = \"467\u00f78=58, 3\" }\n    @{ Old = \"148\u00f76=24, 4\"; New = \"411\u00f77=58, 5\" }\n    @{ Old = \"290\u00f75=58, 0\"; New = \"979\u00f78=122, 3\" }\n    @{ Old = \"340\u00f74=85, 0\"; New = \"760\u00f76=126, 4\" }\n    @{ Old = \"913\u00f76=152, 1\"; New = \"816\u00f74=204, 0\" }\n    @{ Old = \"658\u00f78=82, 2\"; New = \"507\u00f76=84, 3\" }\n    @{ Old = \"925\u00f75=185, 0\"; New = \"816\u00f77=116, 4\" }\n    @{ Old = \"286\u00f76=47, 4\"; New = \"958\u00f79=106, 4\" }\n    @{ Old = \"267\u00f72=133, 1\"; New = \"335\u00f76=55, 5\" }\n    @{ Old = \"902\u00f79=100, 2\"; New = \"341\u00f75=68, 1\" }\n    @{ Old = \"380\u00f78=47, 4\"; New = \"144\u00f75=28, 4\" }\n    @{ Old = \"807\u00f79=89, 6\"; New = \"994\u00f79=110, 4\" }\n    @{ Old = \"294\u00f73=98, 0\"; New = \"188\u00f76=31, 2\" }\n    @{ Old = \"706\u00f72=353, 0\"; New = \"647\u00f72=323, 1\" }\n    @{ Old = \"980\u00f73=326, 2\"; New = \"134\u00f72=67, 0\" }\n    @{ Old = \"732\u00f78=91, 4\"; New = \"372\u00f77=53, 1\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
